# edit.ps1 - applies the "small mod to bc stuff" commit to before.docx
#
# Net effect of the diff (after accounting for pure run-splitting /
# w:proofErr bookkeeping that Word's editor adds automatically as you
# type near a word it doesn't recognise -- it never changes the
# document's visible text, and this headless COM host has no knob to
# force those markers in directly):
#
#   1. In the first paragraph, after "...30 sigma layers in the
#      vertical" insert " (Morrison et al. 2014, Foreman et al. 2014)"
#      right before the following ". Larvae were released...".
#   2. Merge "...randomly within " + "20 x 20 km" into a single run
#      (no text change).
#   3. Merge "...grid cell" + ", if you look ... larvae releas" +
#      "ing from each cell." into a single run and drop the stale
#      _GoBack bookmark that used to sit in the middle of that
#      sentence (no text change).
#   4. Drop the stale lastRenderedPageBreak layout hint in front of
#      "References" (no text change).
#   5. Insert a brand new reference-list paragraph for the Morrison et
#      al. 2014 citation directly after the existing Foreman et al.
#      reference paragraph (and before the North et al. reference
#      paragraph).
#
$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1

# ---------------------------------------------------------------------
# 1) Add the Morrison/Foreman citation right after "...in the vertical"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The ROMS had a ~3 km horizontal resolution and 30 sigma layers in the vertical",
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "The ROMS had a ~3 km horizontal resolution and 30 sigma layers in the vertical (Morrison et al. 2014, Foreman et al. 2014)",
    $wdReplaceOne
) | Out-Null

# ---------------------------------------------------------------------
# 2) Touch the "20 x 20 km" sentence so the two adjacent, identically
#    formatted runs collapse into one (text itself is unchanged).
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "grid cell. Because release occurred randomly within 20 x 20 km",
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "grid cell. Because release occurred randomly within 20 x 20 km",
    $wdReplaceOne
) | Out-Null

# ---------------------------------------------------------------------
# 3) Touch the second "grid cell..." sentence the same way; this also
#    clears out the old _GoBack bookmark that sat inside it.
# ---------------------------------------------------------------------
$apostrophe = [char]0x2019
$gridSentence = "grid cell, if you look at different resolution, won" + $apostrophe + "t have same number of larvae releasing from each cell."
$d.Content.Find.Execute(
    $gridSentence,
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    $gridSentence,
    $wdReplaceOne
) | Out-Null

# ---------------------------------------------------------------------
# 4) Touch "References" so the stale lastRenderedPageBreak hint drops.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "References",
    $false, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "References",
    $wdReplaceOne
) | Out-Null

# ---------------------------------------------------------------------
# 5) Insert the new Morrison et al. reference paragraph after the
#    Foreman et al. reference paragraph, matching its style/formatting.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Foreman, M. G. G.*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
    $newPara = $target.Next()
    $newPara.Range.Text = "Morrison et al. 2014. A model simulation of future oceanic conditions along the british Columbia continental shelf. Part I: Forcing Fields and Initial Conditions. Atmosphere-Ocean, 52, 1-19."
}

Write-Output "done"
